# Crawl-data header fix: the sheet used to list "제품 명" (product name)
# before "브랜드 명" (brand name) in columns C/D. The author re-ordered the
# headers so brand name comes before product name, matching the crawl
# field order described in the commit message:
#   image, brand name, product name, price, detail url
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "브랜드 명"
$ws.Range("D1").Value = "제품 명"

# Leave the selection on the row below the header, as in the saved file.
$null = $ws.Range("A2").Select()
